# Updated cryptos list on Wed Sep 27 14:49:28 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns (D: Price, E: Volume) keep their original text
# formatting (thousand-dot separators, trailing zeros, padding spaces)
# instead of being auto-converted to numbers by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.331.56"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.53%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.604.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.16%  "

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.39%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "213.10"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.77%  "

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.46%  "

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.19%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.38%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.00"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -1.29%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0853"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.60%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.833.14"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.36%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.605.22"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.15%  "

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.01%  "

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -2.53%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.77"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -0.45%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.352.81"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +0.60%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.57"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +6.71%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.0₃0725"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +3.49%  "

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.44%  "

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.90%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.18"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.02%  "

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.57"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.41%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.39%  "

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +1.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.45"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +2.07%  "

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.47%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +1.23%  "

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +0.83%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.443.83"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +8.29%  "

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +1.28%  "

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.34%  "

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.97%  "

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.58%  "

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.94%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.826"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.03%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.84"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +0.79%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.47%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.17"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.932"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -3.00%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.743.38"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +1.29%  "

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.21%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.95"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.09%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "87.60"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +1.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.48"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -2.47%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0500"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.19%  "

$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.35%  "
